$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 24.19710677992798
$ws.Cells.Item(2, 3).Value = 5.768701190967377
$ws.Cells.Item(2, 4).Value = 13.11272690450736
$ws.Cells.Item(2, 5).Value = 12.89901885417056
$ws.Cells.Item(2, 7).Value = 74.98116606619399
$ws.Cells.Item(2, 8).Value = 25.53329367688651
$ws.Cells.Item(2, 10).Value = 8.002234214538847
$ws.Cells.Item(2, 12).Value = 12.83889863573922
$ws.Cells.Item(2, 13).Value = 21.00151155723734

# Row 3
$ws.Cells.Item(3, 2).Value = 24.03659297152049
$ws.Cells.Item(3, 3).Value = 5.399071634615252
$ws.Cells.Item(3, 4).Value = 13.12359722675609
$ws.Cells.Item(3, 5).Value = 12.91542593553118
$ws.Cells.Item(3, 7).Value = 74.53489644095556
$ws.Cells.Item(3, 8).Value = 25.5012434444105
$ws.Cells.Item(3, 10).Value = 7.991239293620393
$ws.Cells.Item(3, 12).Value = 12.85993750338197
$ws.Cells.Item(3, 13).Value = 21.00049311760199

# Row 4
$ws.Cells.Item(4, 2).Value = 23.94490742896246
$ws.Cells.Item(4, 3).Value = 5.15827362218923
$ws.Cells.Item(4, 4).Value = 13.13222070354806
$ws.Cells.Item(4, 5).Value = 12.92619718920119
$ws.Cells.Item(4, 7).Value = 74.27462527788403
$ws.Cells.Item(4, 8).Value = 25.48559282765626
$ws.Cells.Item(4, 10).Value = 7.984259244787004
$ws.Cells.Item(4, 12).Value = 12.87473221582053
$ws.Cells.Item(4, 13).Value = 21.00435956371689

# Row 5
$ws.Cells.Item(5, 2).Value = 23.90930913182861
$ws.Cells.Item(5, 3).Value = 5.05666391490004
$ws.Cells.Item(5, 4).Value = 13.13622469639823
$ws.Cells.Item(5, 5).Value = 12.93076217482463
$ws.Cells.Item(5, 7).Value = 74.17207579934413
$ws.Cells.Item(5, 8).Value = 25.48022837630806
$ws.Cells.Item(5, 10).Value = 7.981357030795324
$ws.Cells.Item(5, 12).Value = 12.88123317578185
$ws.Cells.Item(5, 13).Value = 21.0070659709687

# Row 6
$ws.Cells.Item(6, 2).Value = 23.90350562949735
$ws.Cells.Item(6, 3).Value = 5.039581367455261
$ws.Cells.Item(6, 4).Value = 13.13691913320627
$ws.Cells.Item(6, 5).Value = 12.93153080219508
$ws.Cells.Item(6, 7).Value = 74.15526130152659
$ws.Cells.Item(6, 8).Value = 25.47939881192441
$ws.Cells.Item(6, 10).Value = 7.980871625100495
$ws.Cells.Item(6, 12).Value = 12.88234116366457
$ws.Cells.Item(6, 13).Value = 21.00758369260461

# Row 7
$ws.Cells.Item(7, 2).Value = 23.94442014833068
$ws.Cells.Item(7, 3).Value = 5.156917379068691
$ws.Cells.Item(7, 4).Value = 13.13227271978826
$ws.Cells.Item(7, 5).Value = 12.92625804276062
$ws.Cells.Item(7, 7).Value = 74.27322796314641
$ws.Cells.Item(7, 8).Value = 25.48551637812906
$ws.Cells.Item(7, 10).Value = 7.984220338715878
$ws.Cells.Item(7, 12).Value = 12.87481797890877
$ws.Cells.Item(7, 13).Value = 21.00439148345496

# Row 8
$ws.Cells.Item(8, 2).Value = 24.14035857147844
$ws.Cells.Item(8, 3).Value = 5.644124698772803
$ws.Cells.Item(8, 4).Value = 13.1160703413765
$ws.Cells.Item(8, 5).Value = 12.9045314756711
$ws.Cells.Item(8, 7).Value = 74.82447708463465
$ws.Cells.Item(8, 8).Value = 25.5214070452451
$ws.Cells.Item(8, 10).Value = 7.998490614741055
$ws.Cells.Item(8, 12).Value = 12.84576336657354
$ws.Cells.Item(8, 13).Value = 21.00022907388324

# Row 9
$ws.Cells.Item(9, 2).Value = 24.57735017718058
$ws.Cells.Item(9, 3).Value = 6.489866260906155
$ws.Cells.Item(9, 4).Value = 13.0997711687028
$ws.Cells.Item(9, 5).Value = 12.86744519551905
$ws.Cells.Item(9, 7).Value = 76.01169102866308
$ws.Cells.Item(9, 8).Value = 25.62371578607842
$ws.Cells.Item(9, 10).Value = 8.02466502767207
$ws.Cells.Item(9, 12).Value = 12.80367511915081
$ws.Cells.Item(9, 13).Value = 21.02762848448976

# Row 10
$ws.Cells.Item(10, 2).Value = 24.92802262914608
$ws.Cells.Item(10, 3).Value = 7.044738213595683
$ws.Cells.Item(10, 4).Value = 13.09723774321103
$ws.Cells.Item(10, 5).Value = 12.84354506327334
$ws.Cells.Item(10, 7).Value = 76.94458816511697
$ws.Cells.Item(10, 8).Value = 25.71826697682717
$ws.Cells.Item(10, 10).Value = 8.042804822911041
$ws.Cells.Item(10, 12).Value = 12.78182325191669
$ws.Cells.Item(10, 13).Value = 21.06930239910744

# Row 11
$ws.Cells.Item(11, 2).Value = 25.09339176643983
$ws.Cells.Item(11, 3).Value = 7.28282652241221
$ws.Cells.Item(11, 4).Value = 13.09813508705647
$ws.Cells.Item(11, 5).Value = 12.83339537729845
$ws.Cells.Item(11, 7).Value = 77.38118393314826
$ws.Cells.Item(11, 8).Value = 25.76546049787838
$ws.Cells.Item(11, 10).Value = 8.05082303326634
$ws.Cells.Item(11, 12).Value = 12.7738497592155
$ws.Cells.Item(11, 13).Value = 21.09289811898239

# Row 12
$ws.Cells.Item(12, 2).Value = 25.15680106452997
$ws.Cells.Item(12, 3).Value = 7.37093308997505
$ws.Cells.Item(12, 4).Value = 13.0987693780917
$ws.Cells.Item(12, 5).Value = 12.82965558807069
$ws.Cells.Item(12, 7).Value = 77.54817929754529
$ws.Cells.Item(12, 8).Value = 25.78392904531465
$ws.Cells.Item(12, 10).Value = 8.053825913671881
$ws.Cells.Item(12, 12).Value = 12.77111298228993
$ws.Cells.Item(12, 13).Value = 21.10249565311837

# Row 13
$ws.Cells.Item(13, 2).Value = 25.14311056948463
$ws.Cells.Item(13, 3).Value = 7.352048912156238
$ws.Cells.Item(13, 4).Value = 13.09861968044745
$ws.Cells.Item(13, 5).Value = 12.83045641122644
$ws.Cells.Item(13, 7).Value = 77.51214124144296
$ws.Cells.Item(13, 8).Value = 25.77992501377348
$ws.Cells.Item(13, 10).Value = 8.053180674434447
$ws.Cells.Item(13, 12).Value = 12.77168983142315
$ws.Cells.Item(13, 13).Value = 21.10039927385506

# Row 14
$ws.Cells.Item(14, 2).Value = 25.09859301215774
$ws.Cells.Item(14, 3).Value = 7.290116187795916
$ws.Cells.Item(14, 4).Value = 13.09818137061056
$ws.Cells.Item(14, 5).Value = 12.83308562607687
$ws.Cells.Item(14, 7).Value = 77.39488988438208
$ws.Cells.Item(14, 8).Value = 25.76696796873929
$ws.Cells.Item(14, 10).Value = 8.051070751350851
$ws.Cells.Item(14, 12).Value = 12.7736189405714
$ws.Cells.Item(14, 13).Value = 21.09367446791527

# Row 15
$ws.Cells.Item(15, 2).Value = 25.07142569059981
$ws.Cells.Item(15, 3).Value = 7.251913573798601
$ws.Cells.Item(15, 4).Value = 13.09795123382744
$ws.Cells.Item(15, 5).Value = 12.83470959098441
$ws.Cells.Item(15, 7).Value = 77.32328425743853
$ws.Cells.Item(15, 8).Value = 25.75910907389678
$ws.Cells.Item(15, 10).Value = 8.049774012178778
$ws.Cells.Item(15, 12).Value = 12.77483737179894
$ws.Cells.Item(15, 13).Value = 21.0896414406572

# Row 16
$ws.Cells.Item(16, 2).Value = 24.91732850519333
$ws.Cells.Item(16, 3).Value = 7.028890958034703
$ws.Cells.Item(16, 4).Value = 13.0972203264379
$ws.Cells.Item(16, 5).Value = 12.84422288976637
$ws.Cells.Item(16, 7).Value = 76.9162937588648
$ws.Cells.Item(16, 8).Value = 25.71526659376244
$ws.Cells.Item(16, 10).Value = 8.042276125587708
$ws.Cells.Item(16, 12).Value = 12.78238390043006
$ws.Cells.Item(16, 13).Value = 21.06785331315367

# Row 17
$ws.Cells.Item(17, 2).Value = 24.82425497241596
$ws.Cells.Item(17, 3).Value = 6.888409936568146
$ws.Cells.Item(17, 4).Value = 13.09729678729941
$ws.Cells.Item(17, 5).Value = 12.85024389978049
$ws.Cells.Item(17, 7).Value = 76.66968379211507
$ws.Cells.Item(17, 8).Value = 25.68943908880969
$ws.Cells.Item(17, 10).Value = 8.03761667958209
$ws.Cells.Item(17, 12).Value = 12.78751709831382
$ws.Cells.Item(17, 13).Value = 21.05567183465306

# Row 18
$ws.Cells.Item(18, 2).Value = 24.77127511988043
$ws.Cells.Item(18, 3).Value = 6.806260057885616
$ws.Cells.Item(18, 4).Value = 13.09753372755231
$ws.Cells.Item(18, 5).Value = 12.85377505413886
$ws.Cells.Item(18, 7).Value = 76.52899577834226
$ws.Cells.Item(18, 8).Value = 25.67497767028364
$ws.Cells.Item(18, 10).Value = 8.034914685446513
$ws.Cells.Item(18, 12).Value = 12.79065474243205
$ws.Cells.Item(18, 13).Value = 21.04910223854344

# Row 19
$ws.Cells.Item(19, 2).Value = 24.75343372483684
$ws.Cells.Item(19, 3).Value = 6.778213618895842
$ws.Cells.Item(19, 4).Value = 13.09764709923838
$ws.Cells.Item(19, 5).Value = 12.85498233346841
$ws.Cells.Item(19, 7).Value = 76.48156239771855
$ws.Cells.Item(19, 8).Value = 25.67014904105395
$ws.Cells.Item(19, 10).Value = 8.033996052571277
$ws.Cells.Item(19, 12).Value = 12.79174890353647
$ws.Cells.Item(19, 13).Value = 21.04695304902863

# Row 20
$ws.Cells.Item(20, 2).Value = 24.83410589553498
$ws.Cells.Item(20, 3).Value = 6.903503934240388
$ws.Cells.Item(20, 4).Value = 13.09726867868725
$ws.Cells.Item(20, 5).Value = 12.84959591397972
$ws.Cells.Item(20, 7).Value = 76.69581689781691
$ws.Cells.Item(20, 8).Value = 25.69214773232218
$ws.Cells.Item(20, 10).Value = 8.038114961736387
$ws.Cells.Item(20, 12).Value = 12.78695149842241
$ws.Cells.Item(20, 13).Value = 21.05692339107955

# Row 21
$ws.Cells.Item(21, 2).Value = 25.11164795294538
$ws.Cells.Item(21, 3).Value = 7.308362964880448
$ws.Cells.Item(21, 4).Value = 13.09830212333591
$ws.Cells.Item(21, 5).Value = 12.83231055003013
$ws.Cells.Item(21, 7).Value = 77.42928497455866
$ws.Cells.Item(21, 8).Value = 25.77075758863175
$ws.Cells.Item(21, 10).Value = 8.051691393231868
$ws.Cells.Item(21, 12).Value = 12.77304464667763
$ws.Cells.Item(21, 13).Value = 21.09563177009966

# Row 22
$ws.Cells.Item(22, 2).Value = 25.29760290042686
$ws.Cells.Item(22, 3).Value = 7.561003290439921
$ws.Cells.Item(22, 4).Value = 13.10069381944483
$ws.Cells.Item(22, 5).Value = 12.82161776620019
$ws.Cells.Item(22, 7).Value = 77.9183152123628
$ws.Cells.Item(22, 8).Value = 25.82561341779677
$ws.Cells.Item(22, 10).Value = 8.060369525771922
$ws.Cells.Item(22, 12).Value = 12.76560283958963
$ws.Cells.Item(22, 13).Value = 21.12478805018267

# Row 23
$ws.Cells.Item(23, 2).Value = 25.19795542024892
$ws.Cells.Item(23, 3).Value = 7.427255899641351
$ws.Cells.Item(23, 4).Value = 13.09926040793054
$ws.Cells.Item(23, 5).Value = 12.82726949784202
$ws.Cells.Item(23, 7).Value = 77.65645691583745
$ws.Cells.Item(23, 8).Value = 25.79601889799558
$ws.Cells.Item(23, 10).Value = 8.055755619014866
$ws.Cells.Item(23, 12).Value = 12.76942405531824
$ws.Cells.Item(23, 13).Value = 21.10887544843546

# Row 24
$ws.Cells.Item(24, 2).Value = 24.82965064033877
$ws.Cells.Item(24, 3).Value = 6.896684253895143
$ws.Cells.Item(24, 4).Value = 13.09728078547579
$ws.Cells.Item(24, 5).Value = 12.84988865156802
$ws.Cells.Item(24, 7).Value = 76.68399871461224
$ws.Cells.Item(24, 8).Value = 25.69092194838743
$ws.Cells.Item(24, 10).Value = 8.037889760700926
$ws.Cells.Item(24, 12).Value = 12.78720662520837
$ws.Cells.Item(24, 13).Value = 21.0563562112043

# Row 25
$ws.Cells.Item(25, 2).Value = 24.45373746117221
$ws.Cells.Item(25, 3).Value = 6.272778967503808
$ws.Cells.Item(25, 4).Value = 13.10252265563108
$ws.Cells.Item(25, 5).Value = 12.87688895818002
$ws.Cells.Item(25, 7).Value = 75.67955040661319
$ws.Cells.Item(25, 8).Value = 25.59262610758604
$ws.Cells.Item(25, 10).Value = 8.017776576460728
$ws.Cells.Item(25, 12).Value = 12.81346765113817
$ws.Cells.Item(25, 13).Value = 21.01642077594322
